# Actualización automática 2025-09-18 09:12:30
#
# Adds a new client row ("MAÑAY REAL NOEMI ELIZABETH") between "LUNA PAZMIÑO
# MYRIAM DEL ROCIO" and "MERIZALDE PEREIRA KAREN ELIZABETH" (alphabetical
# order) on the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, and corrects
# the PORCELANATO / septiembre figures for TOSCANO RAMIREZ MONICA CECILIA
# and ZAMBRANO ANGELA MARIA, which ripple into the monthly totals and the
# "CUMPLIMIENTO MENSUAL" roll-up.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------------------
# 1) "VENTAS POR GRUPO": insert the new client row at row 16 (pushes the
#    existing rows 16-25 down to 17-26) and fill it in with zeros.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(16).Insert()

# Use Value2 (rather than Value) so the numeric-format style ("s=2") that
# Insert() already propagated into the blank row survives the write - a
# plain .Value assignment on a still-empty cell resets it to the default
# style.
$ws1.Range("A16").Value2 = "RIOS CARRION ANGEL BENIGNO"
$ws1.Range("B16").Value2 = "MAÑAY REAL NOEMI ELIZABETH"
$ws1.Range("C16:R16").Value2 = 0

# Correct the PORCELANATO value for TOSCANO RAMIREZ MONICA CECILIA (now row 22).
$ws1.Range("M22").Value = -20.74

# Correct the PORCELANATO value for ZAMBRANO ANGELA MARIA (now row 25).
$ws1.Range("M25").Value = 3690.09

# The trailing "x de 23" summary row (now row 26) must read "x de 24" since
# there is one more advisor/client row in the sheet.
$ws1.Range("C26").Value = "0 de 24"
$ws1.Range("D26").Value = "0 de 24"
$ws1.Range("E26").Value = "0 de 24"
$ws1.Range("F26").Value = "0 de 24"
$ws1.Range("G26").Value = "0 de 24"
$ws1.Range("H26").Value = "1 de 24"
$ws1.Range("I26").Value = "0 de 24"
$ws1.Range("J26").Value = "0 de 24"
$ws1.Range("K26").Value = "0 de 24"
$ws1.Range("L26").Value = "1 de 24"
$ws1.Range("M26").Value = "5 de 24"
$ws1.Range("N26").Value = "0 de 24"
$ws1.Range("O26").Value = "0 de 24"
$ws1.Range("P26").Value = "0 de 24"
$ws1.Range("Q26").Value = "0 de 24"
$ws1.Range("R26").Value = "0 de 24"

# ---------------------------------------------------------------------------
# 2) "VENTA MENSUAL": same new row, same two corrected figures (column F is
#    "septiembre" here, the monthly equivalent of "PORCELANATO" above).
# ---------------------------------------------------------------------------
$ws2.Rows.Item(16).Insert()

$ws2.Range("A16").Value2 = "RIOS CARRION ANGEL BENIGNO"
$ws2.Range("B16").Value2 = "MAÑAY REAL NOEMI ELIZABETH"
$ws2.Range("C16:G16").Value2 = 0

# TOSCANO RAMIREZ MONICA CECILIA (now row 22).
$ws2.Range("F22").Value = -20.74

# ZAMBRANO ANGELA MARIA (now row 25).
$ws2.Range("F25").Value = 3690.09

# Monthly total row (now row 26) must reflect the two corrections above.
$ws2.Range("F26").Value = 19884.21

# ---------------------------------------------------------------------------
# 3) "CUMPLIMIENTO MENSUAL": roll-up figures for PORCELANATO (row 12) and the
#    grand TOTAL (row 15) change because VENTA dropped by 6066.55 overall.
# ---------------------------------------------------------------------------
$ws3.Range("D12").Value = 19774.85
$ws3.Range("E12").Value = 23325.2354117774
$ws3.Range("F12").Value = 0.4588123158242369

$ws3.Range("D15").Value = 19884.21
$ws3.Range("E15").Value = 38319.25623249458
$ws3.Range("F15").Value = 0.3416327460734424

# Column D ("VENTA") widens by one character to fit the new figures.
$ws3.Columns.Item(4).ColumnWidth = $ws3.Columns.Item(4).ColumnWidth + 1
